$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.314557909965515
$ws.Range("B1").Value = 1.328503251075745
$ws.Range("C1").Value = 1.427729249000549
$ws.Range("D1").Value = 2.045637369155884
$ws.Range("E1").Value = 4.123985767364502
